$d = $word.ActiveDocument

# 1. Update the "Curso (semestre ideal)" line with the new course/semester values.
$d.Content.Find.Execute(
    "Curso (semestre ideal): EF (9), EM (8), EB (8), EP (10), EQD (9), EQN (10)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Curso (semestre ideal): EF (9), EM (8), EA (6), EP (10), EQD (7), EQN (10)",
    2
) | Out-Null

# 2. Remove the "Requisitos" heading paragraph and the requirement bullet paragraph
#    that follows it (LOB1008 ...), deleting both paragraphs in full.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $text = $p.Range.Text
    if ($text -like "Requisitos*" -or $text -like "LOB1008*") {
        $p.Range.Delete() | Out-Null
    }
}
